$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet lists purchase invoices grouped by vendor. Each vendor group is
# one data row followed (sometimes) by extra "continuation" rows that add
# another invoice for the same vendor, with column F running a cumulative
# total for the group. Three new continuation rows are being added:
#   - under Sr.No 2 (Shree Laxmi Lighting Hub)      -> new row 5
#   - under Sr.No 7 (the "Print ..." vendor)        -> new row 16
#   - under Sr.No 10 (Hatley Technologies)          -> new row 23
# Inserting a real worksheet row at each spot (rather than just overwriting)
# reproduces the row-shift seen throughout the sheet.
# ---------------------------------------------------------------------------

$xlPasteFormats = -4122

function Set-FormatFrom($donorAddr, $targetAddr) {
    $ws.Range($donorAddr).Copy()
    $ws.Range($targetAddr).PasteSpecial($xlPasteFormats)
}

# === 1) New row 5 (continuation of Sr.No 2 / Shree Laxmi Lighting Hub) =====
$ws.Rows("5").Insert()

Set-FormatFrom "A2" "A5"
Set-FormatFrom "B2" "B5"
Set-FormatFrom "C2" "C5"
Set-FormatFrom "D2" "D5"
Set-FormatFrom "E2" "E5"
Set-FormatFrom "F2" "F5"

$ws.Range("B5").Value = 45292
$ws.Range("C5").Value = "SLH/3399"
$ws.Range("D5").Value = "Shree Laxmi Lighting Hub"
$ws.Range("E5").Value = 320
$ws.Range("F5").Formula = "=E4+E5"
Set-FormatFrom "F2" "F5"

# Row 4's own running-total formula moves onto row 5, so F4 goes blank.
$ws.Range("F4").ClearContents()

# === 2) New row 16 (continuation of Sr.No 7 / "Print ...") ================
$ws.Rows("16").Insert()

Set-FormatFrom "A2" "A16"
Set-FormatFrom "B2" "B16"
Set-FormatFrom "C2" "C16"
Set-FormatFrom "D2" "D16"
Set-FormatFrom "E2" "E16"
Set-FormatFrom "F15" "F16"
Set-FormatFrom "C3" "G16"

$ws.Range("B16").Value = 45292
$ws.Range("C16").Value = "2023-24/10627"
$ws.Range("D16").Value = "Print House"
$ws.Range("E16").Value = 6048
$ws.Range("F16").Formula = "=E15+E16"
Set-FormatFrom "F15" "F16"

# Row 15's own running-total formula moves onto row 16, so F15 goes blank;
# its label cell (G15) gets a right/top/bottom box border.
$ws.Range("F15").ClearContents()
Set-FormatFrom "C3" "G15"
$gb = $ws.Range("G15").Borders
$gb.LineStyle = 1
$gb.Weight = 2
$ws.Range("G15").Borders.Item(7).LineStyle = -4142

# === 3) New row 23 (continuation of Sr.No 10 / Hatley Technologies) =======
$ws.Rows("23").Insert()

Set-FormatFrom "A2" "A23"
Set-FormatFrom "B2" "B23"
Set-FormatFrom "C2" "C23"
Set-FormatFrom "D2" "D23"
Set-FormatFrom "E2" "E23"
Set-FormatFrom "F2" "F23"

$ws.Range("B23").Value = 45290
$ws.Range("C23").Value = "PNJ/23-24/3230"
$ws.Range("D23").Value = "Hatley Technologies"
$ws.Range("E23").Formula = "=24019-12974"
$ws.Range("F23").Formula = "=E22+E23"
Set-FormatFrom "F2" "F23"

# Row 22's own running-total formula moves onto row 23, so F22 goes blank.
$ws.Range("F22").ClearContents()

# Recalculate so cached <v> values are correct, then park the selection
# where the user left it.
$wb.Application.Calculate()
$ws.Range("A24").Select()
